# edit.ps1 - applies the surat_tugas_tiket.docx template edit:
#   1. Bumps the embedded OLEObject's internal ObjectID attribute.
#   2. Rewrites the "{{$no}}" placeholder to "${no}" and splits the
#      sentence into several runs, moving the hidden "_GoBack" bookmark
#      to sit right after the new "${no}" placeholder (this also shifts
#      the "_Hlk59436216" bookmark's numeric id, exactly like the diff).
#   3. Collapses the "${" / "tanggal}" runs (and the old "_GoBack" that
#      used to sit between them) back into a single "${tanggal}" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. OLEObject/@ObjectID: _1686057521 -> _1687004657
# ---------------------------------------------------------------------
$oleParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w14:paraId="75716F46" w14:textId="77777777" w:rsidR="00D945DB" w:rsidRPr="00B8635F" w:rsidRDefault="00B94FCC" w:rsidP="00C348A6"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID" w:eastAsia="id-ID"/></w:rPr><w:object w:dxaOrig="1440" w:dyaOrig="1440" w14:anchorId="59516AC9"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_s1398" type="#_x0000_t75" style="position:absolute;left:0;text-align:left;margin-left:12.6pt;margin-top:-55.9pt;width:63.75pt;height:63.75pt;z-index:251657216;visibility:visible;mso-wrap-edited:f" filled="t" fillcolor="#fc0"><v:imagedata r:id="rId6" o:title=""/><w10:wrap type="topAndBottom"/></v:shape><o:OLEObject Type="Embed" ProgID="Word.Picture.8" ShapeID="_x0000_s1398" DrawAspect="Content" ObjectID="_1687004657" r:id="rId7"/></w:object></w:r></w:p>
'@

$oleTable = $d.Tables.Item(1)
$oleCell = $oleTable.Cell(1, 1)
$oleParagraphRange = $oleCell.Range.Paragraphs.Item(1).Range
$oleParagraphRange.InsertXML($oleParagraphXml)

# ---------------------------------------------------------------------
# 2. "Berdasarkan hasil laporan tiket dengan nomor {{$no}} " ->
#    five runs ending "...nomor ${no} " with a relocated _GoBack bookmark
# ---------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Berdasarkan hasil laporan tiket dengan nomor {{`$no}} ")
if (-not $found) {
    throw "Could not find the 'Berdasarkan hasil laporan...' sentence"
}
$sentenceStart = $find.Start
$find.Text = "Berdasarkan hasil laporan tiket dengan nomor `${no} "

function Split-WordRange($rangeStart, $rangeEnd) {
    $piece = $d.Range($rangeStart, $rangeEnd)
    $piece.Font.Bold = 1
    $piece.Font.Bold = 0
}

$b0 = $sentenceStart
$b1 = $sentenceStart + 15   # "Berdasarkan has" | "il laporan tiket dengan nomor "
$b2 = $sentenceStart + 45   # ... | "$"
$b3 = $sentenceStart + 46   # "$" | "{no}"
$b4 = $sentenceStart + 50   # "{no}" | " "   <- _GoBack bookmark goes here
$b5 = $sentenceStart + 51

Split-WordRange $b0 $b1
Split-WordRange $b1 $b2
Split-WordRange $b2 $b3
Split-WordRange $b3 $b4
Split-WordRange $b4 $b5

$goBackRange = $d.Range($b4, $b4)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------
# 3. Merge "${" + "tanggal}" (which used to have _GoBack between them)
#    back into a single "${tanggal}" run.
# ---------------------------------------------------------------------
$tanggalFind = $d.Content
$tanggalFound = $tanggalFind.Find.Execute("Wonosobo, `${tanggal}")
if (-not $tanggalFound) {
    throw "Could not find the 'Wonosobo, `${tanggal}' sentence"
}
$wonosoboStart = $tanggalFind.Start
$placeholderStart = $wonosoboStart + 10
$placeholderEnd = $wonosoboStart + 20
$placeholderRange = $d.Range($placeholderStart, $placeholderEnd)
$placeholderRange.Text = "XXXXXXXXXX"
$placeholderRange2 = $d.Range($placeholderStart, $placeholderStart + 10)
$placeholderRange2.Text = "`${tanggal}"

Write-Host "edit.ps1 completed"
